$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 423, shifting existing rows 423:469 down to 424:470
$ws.Rows.Item(423).Insert()

# Populate the newly inserted row 423 with the new data record
$ws.Range("A423").Value2 = 10
$ws.Range("B423").Value2 = "Vega Modelo de Temuco"
$ws.Range("C423").Value2 = "La Araucanía"
$ws.Range("D423").Value2 = 44858
$ws.Range("E423").Value2 = 9
$ws.Range("F423").Value2 = 100112024
$ws.Range("G423").Value2 = "Choclo"
$ws.Range("H423").Value2 = "Dulce o Americano"
$ws.Range("I423").Value2 = "Primera"
$ws.Range("J423").Value2 = 2500
$ws.Range("K423").Value2 = 600
$ws.Range("L423").Value2 = 600
$ws.Range("M423").Value2 = 600
$ws.Range("N423").Value2 = "$/unidad"
$ws.Range("O423").Value2 = "Argentina"
$ws.Range("P423").Value2 = 600
$ws.Range("Q423").Value2 = 1
$ws.Range("R423").Value2 = "Hortaliza"
